$p = $ppt.ActivePresentation

# --- Slide 1, Title shape: "New Test" -> "PaTESTge" ---
# Title paragraph runs: "New Test" | " " | "1".
# Only the first run's text changes; do a plain text substitution so the
# existing run/paragraph structure (and the other two runs) stay intact.
$slide1 = $p.Slides.Item(1)
$title1 = $slide1.Shapes.Item(1).TextFrame.TextRange
$title1.Replace("New Test", "PaTESTge") | Out-Null

# --- Slide 2, Title shape: "Pare " -> "Pare" + " " (split into two runs) ---
# Title paragraph runs: "Pare " | "4".
# Split the first run into "Pare" and " " while leaving the "4" run alone.
$slide2 = $p.Slides.Item(2)
$title2 = $slide2.Shapes.Item(1).TextFrame.TextRange

$word = $title2.Characters(1, 4)      # "Pare"
$word.Text = "Pare"

$space = $title2.Characters(5, 1)     # " "
$space.Text = " "
